# The deck currently carries two theme parts:
#   - the "live" theme used by the slide master / notes master / handout
#     master (reached through ThemeColorScheme from any of them) is the
#     green "Integral" palette
#   - a second, unused theme part still holds the original "Office Theme"
#     palette
#
# The authored edit swaps the palette back to the standard Office Theme
# colors. Recolor every theme color slot (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) to the Office Theme RGB values via
# ThemeColorScheme.Colors(i).RGB, which is the supported way to edit a
# theme's color scheme through the object model.

$p = $ppt.ActivePresentation

# Office Theme color scheme, in the standard clrScheme slot order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5,
# accent6, hlink, folHlink.
$officeThemeHex = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$tcs = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $officeThemeHex.Count; $i++) {
    $hex = $officeThemeHex[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # OLE/VBA RGB colors are stored little-endian as 0xBBGGRR.
    $oleColor = ($b * 65536) + ($g * 256) + $r
    $tcs.Colors($i).RGB = $oleColor
}
